$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: "Computing Time (sec)" -> "Computing Time (ns)" (shared by C1 and F1)
$ws.Cells.Item(1, 3).Value = "Computing Time (ns)"
$ws.Cells.Item(1, 6).Value = "Computing Time (ns)"

# Row 2 (Instance 6)
$ws.Cells.Item(2, 2).Value = 204
$ws.Cells.Item(2, 3).Value = 91600
$ws.Cells.Item(2, 4).Value = 69.15000000000001
$ws.Cells.Item(2, 6).Value = 97200

# Row 3 (Instance 7)
$ws.Cells.Item(3, 2).Value = 993
$ws.Cells.Item(3, 3).Value = 109700
$ws.Cells.Item(3, 4).Value = 96.97
$ws.Cells.Item(3, 6).Value = 154300

# Row 4 (Instance 8)
$ws.Cells.Item(4, 2).Value = 33
$ws.Cells.Item(4, 3).Value = 61300
$ws.Cells.Item(4, 4).Value = 94.29000000000001
$ws.Cells.Item(4, 6).Value = 72800

# Row 5 (Instance 9)
$ws.Cells.Item(5, 3).Value = 59000
$ws.Cells.Item(5, 6).Value = 64200

# Row 6 (Instance 10)
$ws.Cells.Item(6, 2).Value = 50
$ws.Cells.Item(6, 3).Value = 77700
$ws.Cells.Item(6, 4).Value = 96.15000000000001
$ws.Cells.Item(6, 6).Value = 86400

# Row 7 (Instance 11)
$ws.Cells.Item(7, 2).Value = 91
$ws.Cells.Item(7, 3).Value = 62200
$ws.Cells.Item(7, 4).Value = 85.05
$ws.Cells.Item(7, 6).Value = 67400

# Row 8 (Instance 12)
$ws.Cells.Item(8, 2).Value = 9738
$ws.Cells.Item(8, 3).Value = 117900
$ws.Cells.Item(8, 4).Value = 99.7
$ws.Cells.Item(8, 6).Value = 138200

# Row 9 (Instance 13)
$ws.Cells.Item(9, 2).Value = 106
$ws.Cells.Item(9, 3).Value = 71800
$ws.Cells.Item(9, 4).Value = 81.54000000000001
$ws.Cells.Item(9, 6).Value = 68900

# Row 10 (Instance 14)
$ws.Cells.Item(10, 2).Value = 788
$ws.Cells.Item(10, 3).Value = 142400
$ws.Cells.Item(10, 4).Value = 76.88
$ws.Cells.Item(10, 6).Value = 148800
